$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4721.75
$ws.Range("I51").Value = 3962.3333
$ws.Range("K51").Value = 3962.3333
$ws.Range("M51").Value = -3478.3333
$ws.Range("H100").Value = 1271
$ws.Range("I100").Value = 1271
$ws.Range("K100").Value = 1271
$ws.Range("M100").Value = -730
$ws.Range("H113").Value = 20116.166
$ws.Range("I113").Value = 21626.727
$ws.Range("K113").Value = 21626.727
$ws.Range("M113").Value = -18372.727
$ws.Range("H116").Value = 19262.5
$ws.Range("J116").Value = 7728.5713
$ws.Range("L116").Value = 7728.5713
$ws.Range("N116").Value = -14612.5713
$ws.Range("H132").Value = 965.6111
$ws.Range("I132").Value = 962.51514
$ws.Range("K132").Value = 2887.54542
$ws.Range("M132").Value = -357.5454199999999
$ws.Range("H135").Value = 774.1539
$ws.Range("I135").Value = 797
$ws.Range("K135").Value = 7173
$ws.Range("M135").Value = -4638
$ws.Range("H139").Value = 47299.668
$ws.Range("J139").Value = 47299.668
$ws.Range("L139").Value = 47299.668
$ws.Range("N139").Value = -57579.668
$ws.Range("H140").Value = 53245.6
$ws.Range("J140").Value = 53245.6
$ws.Range("L140").Value = 53245.6
$ws.Range("N140").Value = -63605.6

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4039.827
$ws.Range("I32").Value = 2151.8718
$ws.Range("K32").Value = 2151.8718
$ws.Range("M32").Value = -1864.8718

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 16199.8
$ws.Range("I97").Value = 7999.3335
$ws.Range("J97").Value = 19714.285
$ws.Range("K97").Value = 7999.3335
$ws.Range("L97").Value = 19714.285
$ws.Range("M97").Value = -7008.3335
$ws.Range("N97").Value = -21696.285
$ws.Range("H99").Value = 999
$ws.Range("I99").Value = 999
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 999
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 499
$ws.Range("N99").Value = ""
$ws.Range("H134").Value = 9974.270500000001
$ws.Range("I134").Value = 10204.566
$ws.Range("J134").Value = 8987.286
$ws.Range("K134").Value = 30613.698
$ws.Range("L134").Value = 26961.858
$ws.Range("M134").Value = -28078.698
$ws.Range("N134").Value = -32031.858

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2657.524
$ws.Range("I31").Value = 1517.2941
$ws.Range("J31").Value = 7503.5
$ws.Range("K31").Value = 1517.2941
$ws.Range("L31").Value = 7503.5
$ws.Range("M31").Value = -1222.2941
$ws.Range("N31").Value = -8093.5
$ws.Range("H34").Value = 2657.524
$ws.Range("I34").Value = 1517.2941
$ws.Range("J34").Value = 7503.5
$ws.Range("K34").Value = 1517.2941
$ws.Range("L34").Value = 7503.5
$ws.Range("M34").Value = -1315.2941
$ws.Range("N34").Value = -7907.5
$ws.Range("H51").Value = 29653.666
$ws.Range("I51").Value = 22222
$ws.Range("J51").Value = 31140
$ws.Range("K51").Value = 22222
$ws.Range("L51").Value = 31140
$ws.Range("M51").Value = -21486
$ws.Range("N51").Value = -32612
$ws.Range("H58").Value = 1891736.1
$ws.Range("I58").Value = 3345563.2
$ws.Range("J58").Value = 1760.8
$ws.Range("K58").Value = 3345563.2
$ws.Range("L58").Value = 1760.8
$ws.Range("M58").Value = -3345360.2
$ws.Range("N58").Value = -2166.8
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("H60").Value = 11098.261
$ws.Range("J60").Value = 11098.261
$ws.Range("L60").Value = 11098.261
$ws.Range("N60").Value = -12120.261
$ws.Range("H61").Value = 29653.666
$ws.Range("I61").Value = 22222
$ws.Range("J61").Value = 31140
$ws.Range("K61").Value = 22222
$ws.Range("L61").Value = 31140
$ws.Range("M61").Value = -21874
$ws.Range("N61").Value = -31836
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31748
$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -98736
$ws.Range("H132").Value = 1515.2122
$ws.Range("I132").Value = 1032.375
$ws.Range("K132").Value = 3097.125
$ws.Range("M132").Value = -567.125
$ws.Range("H134").Value = 2175.5
$ws.Range("I134").Value = 1347.7059
$ws.Range("J134").Value = 6866.3335
$ws.Range("K134").Value = 4043.1177
$ws.Range("L134").Value = 20599.0005
$ws.Range("M134").Value = -1508.1177
$ws.Range("N134").Value = -25669.0005
$ws.Range("H136").Value = 1891736.1
$ws.Range("I136").Value = 3345563.2
$ws.Range("J136").Value = 1760.8
$ws.Range("K136").Value = 10036689.6
$ws.Range("L136").Value = 5282.4
$ws.Range("M136").Value = -10034139.6
$ws.Range("N136").Value = -10382.4

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 265.42856
$ws.Range("I40").Value = 226
$ws.Range("J40").Value = 295
$ws.Range("K40").Value = 904
$ws.Range("L40").Value = 1180
$ws.Range("M40").Value = -835
$ws.Range("N40").Value = -1318
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = ""

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5799.75
$ws.Range("I70").Value = 6244.222
$ws.Range("K70").Value = 6244.222
$ws.Range("M70").Value = -5974.222
$ws.Range("H73").Value = 5799.75
$ws.Range("I73").Value = 6244.222
$ws.Range("K73").Value = 6244.222
$ws.Range("M73").Value = -5308.222
$ws.Range("H132").Value = 1285157
$ws.Range("I132").Value = 1749953.1
$ws.Range("K132").Value = 5249859.300000001
$ws.Range("M132").Value = -5247329.300000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 370.41666
$ws.Range("J55").Value = 441.8
$ws.Range("L55").Value = 441.8
$ws.Range("N55").Value = -787.8
$ws.Range("H132").Value = 1624.8788
$ws.Range("I132").Value = 1220.3334
$ws.Range("K132").Value = 3661.0002
$ws.Range("M132").Value = -1131.0002
$ws.Range("H136").Value = 3944.762
$ws.Range("I136").Value = 2099.5557
$ws.Range("J136").Value = 5328.6665
$ws.Range("K136").Value = 6298.6671
$ws.Range("L136").Value = 15985.9995
$ws.Range("M136").Value = -3748.6671
$ws.Range("N136").Value = -21085.9995

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
$ws.Range("H123").Value = 47536.09
$ws.Range("J123").Value = 47536.09
$ws.Range("L123").Value = 47536.09
$ws.Range("N123").Value = -57336.09
$ws.Range("H132").Value = 5675.9473
$ws.Range("I132").Value = 1127.4667
$ws.Range("J132").Value = 8642.348
$ws.Range("K132").Value = 3382.4001
$ws.Range("L132").Value = 25927.044
$ws.Range("M132").Value = -852.4000999999998
$ws.Range("N132").Value = -30987.044
$ws.Range("H136").Value = 25255464
$ws.Range("I136").Value = 34725110
$ws.Range("K136").Value = 104175330
$ws.Range("M136").Value = -104172780

